# Delete (Teran et al., 2019) data from VEGF:NRP1
#
# The VEGFA165_NRP1 worksheet has two rows (Teran et al., 2019 (chimera) and
# Teran et al., 2019 (monomer)) that need to be removed. Deleting the entire
# rows shifts the rows below them up and Excel automatically drops the now
# unused shared-string entries for those two references.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_NRP1")

# Rows 6 and 7 hold the "Teran et al., 2019 (chimera)" and
# "Teran et al., 2019 (monomer)" entries - remove them, shifting rows 8-9 up.
$ws.Range("A6:A7").EntireRow.Delete()

# Make the VEGFA165_NRP1 sheet the active one and select the rows that used
# to be the last two entries (now rows 6:7) to match the author's final view.
$ws.Activate()
$ws.Range("A6:D7").Select()
